$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (data row 1): convert "Numero" P1 text into numeric 1 ---
$ws.Range("A2").Value = 1

# --- Row 3 (new data row 2): Numero = 2 ---
$ws.Range("A3").Value = 2

# B3: "A FAIRE " (already has style fontId5/vertical-top from source file)
$ws.Range("B3").Value = 'A FAIRE '

# D3 and E3 share a new wrapped style based on B3's font (Calibri12) + wrap text
$ws.Range("B3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").WrapText = $true
$ws.Range("D3").Value = 'La balise meta "keywords" était utilisée autrefois pour indiquer aux moteurs de recherche les mots-clés associés à une page web. Cependant, cette balise est devenue obsolète et n''est plus prise en compte par les moteurs de recherche, notamment Google, qui ne l''utilise plus depuis 2009. Il est donc inutile de l''utiliser dans votre code HTML.'

$ws.Range("B3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("E3").WrapText = $true
$ws.Range("E3").Value = 'pour améliorer la référence du site on va plutôt utiliser d''autres éléments tels que le contenu de qualité, la structure du site et les balises meta "description" et "title'

# G3: hyperlink cell, styled like G2 (Lien hypertexte style)
$ws.Range("G3").Value = 'https://www.balisemeta.com/raison-du-declin-balise-keywords.html'
$ws.Hyperlinks.Add($ws.Range("G3"), 'https://www.balisemeta.com/raison-du-declin-balise-keywords.html')
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)

# C3: rich text cell styled like D2/E2 (wrap, left/top aligned issue-description look)
$ws.Range("D2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 'La balise meta "keywords" est obsolète et ne devrait pas être utilisée(2récurence sur la page Index.html et Page2.html) '
$c3run = $ws.Range("C3").Characters(71, 50)
$c3run.Font.Size = 11
$c3run.Font.Color = 255
$c3run.Font.Name = "Arial"

# H1: new "Autre" header cell, bold Calibri on theme accent fill
$ws.Range("B2").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = 'Autre'
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Interior.ThemeColor = 10
$ws.Range("H1").Interior.Pattern = 1

# H3: same wrapped style as D3/E3
$ws.Range("B3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("H3").WrapText = $true
$ws.Range("H3").Value = 'remplacer par ?<meta name="robots" content="index, follow">  '

# F3: "Action recommandee" style like F2
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = 'compléter les balise (qui sont mal utlisé dans le site) en utilisant mieux les balises "description" et "title".Suppression de la ligne sur les 2 pages'

# Row height / column width / selection cosmetic tweaks
$ws.Rows.Item(3).RowHeight = 120
$ws.Columns.Item(5).ColumnWidth = 24.1
$ws.Range("F3").Select()

Write-Host "edit complete"
